$wb = $excel.ActiveWorkbook

$alc = $wb.Worksheets.Item("ALC")
$arm = $wb.Worksheets.Item("ARM")
$bsm = $wb.Worksheets.Item("BSM")
$crp = $wb.Worksheets.Item("CRP")
$cul = $wb.Worksheets.Item("CUL")
$gsm = $wb.Worksheets.Item("GSM")
$ltw = $wb.Worksheets.Item("LTW")
$wvr = $wb.Worksheets.Item("WVR")

$alc.Range("H33").Value = 819.9643
$alc.Range("I33").Value = 547.8261
$alc.Range("K33").Value = 547.8261
$alc.Range("M33").Value = -318.8261

$alc.Range("H113").Value = 4999.5
$alc.Range("I113").Value = 4999
$alc.Range("K113").Value = 4999
$alc.Range("M113").Value = -1745

$alc.Range("H127").Value = 6370.75
$alc.Range("I127").Value = 6370.75
$alc.Range("K127").Value = 19112.25
$alc.Range("M127").Value = -14152.25

$alc.Range("H129").Value = 4053.1428
$alc.Range("I129").Value = 4053.1428
$alc.Range("K129").Value = 12159.4284
$alc.Range("M129").Value = -7159.428400000001

$alc.Range("H137").Value = 1610.8462
$alc.Range("I137").Value = 1254.1
$alc.Range("K137").Value = 3762.3
$alc.Range("M137").Value = -1212.3

$alc.Range("H138").Value = 2323
$alc.Range("I138").Value = 2414.4285
$alc.Range("J138").Value = 1683
$alc.Range("K138").Value = 7243.2855
$alc.Range("L138").Value = 5049
$alc.Range("M138").Value = -2103.2855
$alc.Range("N138").Value = -15329

$alc.Range("H141").Value = 978.7895
$alc.Range("I141").Value = 982.17645
$alc.Range("K141").Value = 2946.52935
$alc.Range("M141").Value = 2233.47065

$arm.Range("H64").Value = 0
$arm.Range("J64").Value = 0
$arm.Range("L64").Value = 0
$arm.Range("N64").ClearContents()

$arm.Range("H67").Value = 0
$arm.Range("J67").Value = 0
$arm.Range("L67").Value = 0
$arm.Range("N67").ClearContents()

$arm.Range("H74").Value = 1458.4
$arm.Range("I74").Value = 1458.4
$arm.Range("K74").Value = 1458.4
$arm.Range("M74").Value = -584.4000000000001

$arm.Range("H77").Value = 1458.4
$arm.Range("I77").Value = 1458.4
$arm.Range("K77").Value = 7292
$arm.Range("M77").Value = -2924

$arm.Range("H97").Value = 729.4231
$arm.Range("I97").Value = 487.65216
$arm.Range("J97").Value = 2583
$arm.Range("K97").Value = 487.65216
$arm.Range("L97").Value = 2583
$arm.Range("M97").Value = 8.347840000000019
$arm.Range("N97").Value = -3575

$bsm.Range("H62").Value = 48777.5
$bsm.Range("J62").Value = 48777.5
$bsm.Range("L62").Value = 48777.5
$bsm.Range("N62").Value = -50149.5

$bsm.Range("H65").Value = 48777.5
$bsm.Range("J65").Value = 48777.5
$bsm.Range("L65").Value = 146332.5
$bsm.Range("N65").Value = -153196.5

$bsm.Range("H97").Value = 13644.333
$bsm.Range("I97").Value = 10399.857
$bsm.Range("K97").Value = 10399.857
$bsm.Range("M97").Value = -9408.857

$crp.Range("H132").Value = 2532.75
$crp.Range("I132").Value = 2781.8096
$crp.Range("K132").Value = 8345.4288
$crp.Range("M132").Value = -5815.4288

$cul.Range("H4").Value = 250000460
$cul.Range("I4").Value = 200000460
$cul.Range("K4").Value = 600001380
$cul.Range("M4").Value = -600001268

$cul.Range("H97").Value = 737.25
$cul.Range("I97").Value = 450
$cul.Range("K97").Value = 1350
$cul.Range("M97").Value = -854

$cul.Range("H98").Value = 487
$cul.Range("I98").Value = 399
$cul.Range("J98").Value = 504.6
$cul.Range("K98").Value = 1197
$cul.Range("L98").Value = 1513.8
$cul.Range("N98").Value = -4509.8
$cul.Range("M98").Value = 301

$cul.Range("H113").Value = 226.25
$cul.Range("I113").Value = 201.5
$cul.Range("J113").Value = 251
$cul.Range("K113").Value = 604.5
$cul.Range("L113").Value = 753
$cul.Range("M113").Value = 1565.5
$cul.Range("N113").Value = -5093

$cul.Range("H125").Value = 0
$cul.Range("I125").Value = 0
$cul.Range("J125").Value = 0
$cul.Range("K125").Value = 0
$cul.Range("L125").Value = 0
$cul.Range("M125").ClearContents()
$cul.Range("N125").ClearContents()

$cul.Range("H129").Value = 2326.6365
$cul.Range("J129").Value = 2913.2856
$cul.Range("L129").Value = 8739.856800000001
$cul.Range("N129").Value = -18739.8568

$cul.Range("H131").Value = 1052.579
$cul.Range("I131").Value = 1100
$cul.Range("J131").Value = 1049.9445
$cul.Range("K131").Value = 3300
$cul.Range("L131").Value = 3149.8335
$cul.Range("M131").Value = 1740
$cul.Range("N131").Value = -13229.8335

$gsm.Range("H7").Value = 1234
$gsm.Range("I7").Value = 1234
$gsm.Range("K7").Value = 1234
$gsm.Range("M7").Value = -1122

$gsm.Range("H8").Value = 1234
$gsm.Range("I8").Value = 1234
$gsm.Range("K8").Value = 1234
$gsm.Range("M8").Value = -1095

$gsm.Range("H33").Value = 24999
$gsm.Range("J33").Value = 24999
$gsm.Range("L33").Value = 24999
$gsm.Range("N33").Value = -25503

$gsm.Range("H43").Value = 13384.5
$gsm.Range("I43").Value = 4267
$gsm.Range("J43").Value = 22502
$gsm.Range("K43").Value = 4267
$gsm.Range("L43").Value = 22502
$gsm.Range("M43").Value = -4116
$gsm.Range("N43").Value = -22804

$gsm.Range("H80").Value = 2708.8
$gsm.Range("I80").Value = 2399.5
$gsm.Range("J80").Value = 2915
$gsm.Range("K80").Value = 2399.5
$gsm.Range("L80").Value = 2915
$gsm.Range("M80").Value = -1401.5
$gsm.Range("N80").Value = -4911

$gsm.Range("H83").Value = 2708.8
$gsm.Range("I83").Value = 2399.5
$gsm.Range("J83").Value = 2915
$gsm.Range("K83").Value = 11997.5
$gsm.Range("L83").Value = 14575
$gsm.Range("M83").Value = -7005.5
$gsm.Range("N83").Value = -24559

$ltw.Range("H23").Value = 9000
$ltw.Range("I23").Value = 9000
$ltw.Range("K23").Value = 9000
$ltw.Range("M23").Value = -8770

$ltw.Range("H43").Value = 20004.334
$ltw.Range("J43").Value = 20004.334
$ltw.Range("L43").Value = 20004.334
$ltw.Range("N43").Value = -20390.334

$ltw.Range("H46").Value = 1678.1305
$ltw.Range("J46").Value = 2237.25
$ltw.Range("L46").Value = 2237.25
$ltw.Range("N46").Value = -2613.25

$ltw.Range("H132").Value = 2820
$ltw.Range("I132").Value = 2801.5715
$ltw.Range("K132").Value = 8404.7145
$ltw.Range("M132").Value = -5874.7145

$ltw.Range("H136").Value = 4440.1333
$ltw.Range("J136").Value = 6000.6665
$ltw.Range("L136").Value = 18001.9995
$ltw.Range("N136").Value = -23101.9995

$wvr.Range("H136").Value = 2447.9443
$wvr.Range("I136").Value = 1397.9286
$wvr.Range("J136").Value = 6123
$wvr.Range("K136").Value = 4193.7858
$wvr.Range("L136").Value = 18369
$wvr.Range("M136").Value = -1643.7858
$wvr.Range("N136").Value = -23469
